$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.354.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.471.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.17%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -1.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.469.19"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.138"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.165"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.338"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.955.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000175"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.211.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.491.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.55%  "
$ws.Range("E19").Value = "  +3.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "366.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.18%  "
$ws.Range("E23").Value = "  -2.76%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("E26").Value = "  -5.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.995"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.622.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0955"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "528.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("E33").Value = "  -3.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("E36").Value = "  -3.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.63%  "
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("E41").Value = "  -2.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "144.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0274"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("E49").Value = "  -3.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0746"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.97%  "
